$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.814.58"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.17%  "
$ws.Range("D3").Value = "'3.486.87"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.85%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'568.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.95%  "
$ws.Range("D6").Value = "'182.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.32%  "
$ws.Range("D7").Value = "'0.613"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.46%  "
$ws.Range("D8").Value = "'3.480.85"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.92%  "
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("E10").Value = "  +1.90%  "
$ws.Range("D12").Value = "'53.51"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.96%  "
$ws.Range("E13").Value = "  -1.81%  "
$ws.Range("D14").Value = "'9.37"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.89%  "
$ws.Range("D15").Value = "'4.048.45"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.87%  "
$ws.Range("D16").Value = "'19.13"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.19%  "
$ws.Range("D17").Value = "'68.724.57"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.11%  "
$ws.Range("D18").Value = "'3.488.25"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.80%  "
$ws.Range("D19").Value = "'12.25"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.16%  "
$ws.Range("E20").Value = "  -1.53%  "
$ws.Range("D21").Value = "'538.25"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +12.53%  "
$ws.Range("E22").Value = "  -3.70%  "
$ws.Range("D23").Value = "'19.32"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.68%  "
$ws.Range("D24").Value = "'4.97"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.77%  "
$ws.Range("D25").Value = "'4.37"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("D26").Value = "'94.13"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.20%  "
$ws.Range("D27").Value = "'2.89"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.04%  "
$ws.Range("D28").Value = "'10.71"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.18%  "
$ws.Range("D29").Value = "'8.98"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.49%  "
$ws.Range("D30").Value = "'31.35"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.27%  "
$ws.Range("D31").Value = "'7.16"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.91%  "
$ws.Range("D32").Value = "'12.49"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.89%  "
$ws.Range("E33").Value = "  -3.65%  "
$ws.Range("E34").Value = "  -5.80%  "
$ws.Range("D35").Value = "'571.40"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.91%  "
$ws.Range("E36").Value = "  -0.08%  "
$ws.Range("D37").Value = "'37.64"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.15%  "
$ws.Range("D38").Value = "'3.01"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.25%  "
$ws.Range("D39").Value = "'0.394"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.18%  "
$ws.Range("D40").Value = "'0.0₃0759"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.34%  "
$ws.Range("D41").Value = "'3.09"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.41%  "
$ws.Range("E42").Value = "  -5.29%  "
$ws.Range("D43").Value = "'3.30"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.69%  "
$ws.Range("D44").Value = "'3.209.36"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.07%  "
$ws.Range("E45").Value = "  -5.36%  "
$ws.Range("D46").Value = "'3.42"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.86%  "
$ws.Range("D47").Value = "'0.0434"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.76%  "
$ws.Range("D48").Value = "'9.01"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.98%  "
$ws.Range("D49").Value = "'0.133"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.39%  "
$ws.Range("D50").Value = "'0.999"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.08%  "
$ws.Range("D51").Value = "'136.56"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.25%  "
